# Applies the review database update:
#  - G21 changes from "yes" to "no"
#  - G22 changes from "yes" to "confirm"
#  - Active cell / selection moves from G20 to G23

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G21").Value = "no"
$ws.Range("G22").Value = "confirm"

$ws.Range("G23").Select()
